$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.01744091126476769
$ws.Range("C2").Value = 0.2093036482578658

$ws.Range("B3").Value = 0.09767403987396309
$ws.Range("C3").Value = 0.2875275565943013

$ws.Range("B4").Value = 0.4739433067654015
$ws.Range("C4").Value = 0.09561128821532201

$ws.Range("B5").Value = 0.9566838304586693
$ws.Range("C5").Value = 0.3714701371556398

$ws.Range("B6").Value = 0.7269852369805344
$ws.Range("C6").Value = 0.4200655772301339

$ws.Range("B7").Value = 0.4954537119737896
$ws.Range("C7").Value = 0.03902275302677104

$ws.Range("B8").Value = 0.01023569464683533
$ws.Range("C8").Value = 0.2692729568481445

$ws.Range("B9").Value = 0.07346130622333351
$ws.Range("C9").Value = 0.215916035453129

$ws.Range("B10").Value = 0.7337921039974529
$ws.Range("C10").Value = 0.5388144169382896
